$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 corresponds to Requirement 12 (data row under header row 1).
# Completeness score changes from 1 to 0, and a Completeness-Reasons
# value of "1.5.3" is now provided.
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = "1.5.3"
